$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$handedBackStatus = "Handed back: in sync with en-US"

# --- Status column updates (all four cells share the same underlying text) ---
$wsOverview.Range("E2").Value = $handedBackStatus
$wsOverview.Range("F2").Value = $handedBackStatus
$wsZhCn.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("C2").Value = $handedBackStatus

# --- zh-cn sheet: handback datetime + cleared error detail ---
$wsZhCn.Range("K2").Value = "2016-08-16 06:43:29"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: handback datetime + cleared error detail ---
$wsDeDe.Range("K2").Value = "2016-08-16 06:43:36"
$wsDeDe.Range("P2").Value = ""

# --- Column width updates ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
